$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos value changed to the professor identifier string ---
$ws.Range("B10").Value = "5840938 - Marcelo Rodrigues de Holanda"
$ws.Range("C10").Value = "5840938 - Marcelo Rodrigues de Holanda"

# --- Row 13: now "Programa resumido:" / "Semestral" (used to be blank label / professor) ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: now "Short syllabus:" only (used to hold the short program text) ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# --- Row 15: now "Programa:" / "01/01/2012" (copy the existing text cell so it
#     is stored as shared text, not auto-converted to a date serial) ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: now "Syllabus:" only (used to hold the long program description) ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()

# --- Row 17: now "Avaliação:" only, normal height (used to be "Syllabus:" tall row) ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# --- Row 18: now "Método:" / professor id (used to be "Avaliação:" only) ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B10").Copy($ws.Range("B18"))
$ws.Range("C10").Copy($ws.Range("C18"))
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: label renamed "Método:" -> "Critério:" (values stay as-is) ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20: label renamed "Critério:" -> "Norma de recuperação:" ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: label renamed "Norma de recuperação:" -> "Bibliografia:", taller row ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22 (old Bibliografia: + long bibliography text) removed entirely ---
$ws.Rows.Item(22).Delete()

Write-Host "Edit complete"
